# Apply "Add bio info, images and sample data for bio-master d" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bio_master_D")

# --- New data values in columns B (Common Name), E (Bio Generated) and
#     H (Bio-verified) for the first several species rows. ---

# Row 2
$ws.Range("E2").Value = "N"

# Row 3
$ws.Range("B3").Value = "Kuranda Tree Frog"
$ws.Range("E3").Value = "Y"
$ws.Range("H3").Value = "Common Name?"

# Row 4
$ws.Range("B4").Value = "Peron's Tree Frog"
$ws.Range("H4").Value = "Common Name?"

# Row 5
$ws.Range("H5").Value = "Y"

# Row 6
$ws.Range("H6").Value = "Y"

# Row 7
$ws.Range("B7").Value = "Orange Thighed Tree Frog"
$ws.Range("H7").Value = "Common Name?"

# Row 8
$ws.Range("H8").Value = "Y"

# Row 9
$ws.Range("H9").Value = "Y"

# Row 10
$ws.Range("H10").Value = "Y"

# Row 11
$ws.Range("B11").Value = "Purple-crowned fairy wren"
$ws.Range("H11").Value = "Common Name?"

# --- Column H width widened to fit the new "Common Name?" values ---
$ws.Columns.Item(8).ColumnWidth = 17.08984375

# --- View state: frozen pane scrolled up and selection moved ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B12").Select()

# --- Application window geometry ---
$excel.ActiveWindow.WindowState = -4143
$excel.Width = 25600
$excel.Height = 15460
